$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 458; this shifts rows 458:503 down to 459:504
# preserving all of their existing values/formatting.
$ws.Rows.Item(458).Insert()

# Populate the newly inserted row 458 with this week's record (same
# dimension values as the former row 458, new measurement columns).
$ws.Range("A458").Value = 10
$ws.Range("B458").Value = "Vega Modelo de Temuco"
$ws.Range("C458").Value = "La Araucanía"
$ws.Range("D458").Value = 45194
$ws.Range("E458").Value = 9
$ws.Range("F458").Value = "Fruta"
$ws.Range("G458").Value = 100102
$ws.Range("H458").Value = "Cítricos"
$ws.Range("I458").Value = 100102006
$ws.Range("J458").Value = "Pomelo"
$ws.Range("K458").Value = "Start Ruby"
$ws.Range("L458").Value = "Primera"
$ws.Range("M458").Value = 125
$ws.Range("N458").Value = 15000
$ws.Range("O458").Value = 15000
$ws.Range("P458").Value = 15000
$ws.Range("Q458").Value = "$/bandeja 15 kilos granel"
$ws.Range("R458").Value = "Región de O'Higgins"
$ws.Range("S458").Value = 1000
$ws.Range("T458").Value = 15
